# feat: add 2022-Q4 data
#
# - Insert a new "2022-Q4" worksheet (a copy of the existing "2022-Q3"
#   sheet, so headers/layout/styling match) positioned right before the
#   current "2022-Q3" sheet, then update its fund figures.
# - Update the "总计" (summary) sheet: insert a new row for "2022-Q4"
#   right after the header row, shifting the existing "2022-Q3" /
#   "2021-Q3" rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet to create the new "2022-Q4" sheet,
#    inserted immediately before the original "2022-Q3" sheet.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Update the quarterly fund figures for 2022-Q4 (keep them as text, same
# as the rest of the column, by forcing a text format before assigning,
# then copying the plain formatting back from a neighbouring cell so no
# stray number format lingers on the cells).
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.36"
$q4.Range("E2").Value = "96.22"
$q4.Range("F2").Value = "0.94"
$q4.Range("G2").Value = "0.0034"
$q4.Range("H2").Value = 5

$q4.Range("B2").Copy() | Out-Null
$q4.Range("D2:G2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet with a new row for 2022-Q4.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()

# Copy formatting from the row below (the old "2022-Q3" row) so the new
# row's styles match the rest of the table.
$zj.Range("A3:D3").Copy() | Out-Null
$zj.Range("A2:D2").PasteSpecial(-4122) | Out-Null

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0

# Renumber the index column for the rows that shifted down.
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2

Write-Output "2022-Q4 sheet added"
